$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = 1556
$ws.Range("E2").Value = 46200502250
$ws.Range("X2").Value = "DN4127450128761"

# Row 3 updates
$ws.Range("A3").Value = 1557
$ws.Range("E3").Value = 46200502251
$ws.Range("X3").Value = "DN4127450128762"
